$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Recordings" column header (H1), bold like the rest of row 1
$ws.Range("H1").Value = "Recordings"
$ws.Range("H1").Font.Bold = $true

# Row 5 - SQL Saturday Baton Rouge 2023 - BI Edition: registration numbers
# (write E5 before D5 so the shared "no show rate" formula in F5 recalculates
# against both updated precedents)
$ws.Range("E5").Value = 106
$ws.Range("D5").Value = 236

# Row 6 - SQL Saturday Costa Rica 2023: attended count + recordings count
$ws.Range("E6").Value = 280
$ws.Range("H6").Value = 640

# New event rows
$ws.Range("A13").Value = "SQL Saturday Boston 2023 "
$ws.Range("A14").Value = "SQL Saturday Columbus"
$ws.Range("A15").Value = "SQL Saturday SLC"

# Leave selection where the author left it
$ws.Range("A13").Select()
